$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated source counts (C/D/E) and recomputed statistics (F delta, G AAPC, H IC_Inf, I IC_Sup)
# Note: values below are written into columns C-F and the *old* H/I/J first;
# deleting the now-empty column G afterwards shifts H->G, I->H, J->I to match
# the target layout (AAPC inserted, trailing IC_Sup column dropped).
$newData = @(
    @{row=2;  C=7669083; D=6755553; E=7189437; F=-6.25428098770088;  AAPC=-3.27450924051266; ICInf=-3.32445223087505; ICSup=-3.22454044939358},
    @{row=3;  C=452208;  D=538188;  E=702538;  F=55.3572692212433;   AAPC=25.0627851107555;  ICInf=24.8293872602317;  ICSup=25.2966193533646},
    @{row=4;  C=7216875; D=6217365; E=6486899; F=-10.1148488784966;  AAPC=-5.35078757963664; ICInf=-5.40171656335717; ICSup=-5.29983117722215},
    @{row=5;  C=4133201; D=3652373; E=3767336; F=-8.85185598280848;  AAPC=-4.64094032801387; ICInf=-4.70830026528327; ICSup=-4.57353277525264},
    @{row=6;  C=102905;  D=117038;  E=176370;  F=71.3910888683737;   AAPC=32.5407951070158;  ICInf=32.0263702388432;  ICSup=33.0572243697993},
    @{row=7;  C=4030296; D=3535335; E=3590966; F=-10.9006881876666;  AAPC=-5.73894993642492; ICInf=-5.80672704613122; ICSup=-5.67112405745457},
    @{row=8;  C=1010395; D=874689;  E=996639;  F=-1.36144775063218;  AAPC=-0.713478646382337;ICInf=-0.85377824830708; ICSup=-0.572980509627385},
    @{row=9;  C=80055;   D=94798;   E=137368;  F=71.5920304790457;   AAPC=32.1680744891262;  ICInf=31.5904945840354;  ICSup=32.7481895206928},
    @{row=10; C=930340;  D=779891;  E=859271;  F=-7.63903519143539;  AAPC=-4.06505764993311; ICInf=-4.20867527295324; ICSup=-3.92122470447432},
    @{row=11; C=1060737; D=934158;  E=1033925; F=-2.52767651170837;  AAPC=-1.31910539884393; ICInf=-1.45512757595543; ICSup=-1.18289546936936},
    @{row=12; C=110178;  D=131747;  E=165242;  F=49.9773094447167;   AAPC=22.6620764704427;  ICInf=22.1967028251339;  ICSup=23.1292224436845},
    @{row=13; C=950559;  D=802411;  E=868683;  F=-8.61345797578057;  AAPC=-4.57820376772015; ICInf=-4.71964304357978; ICSup=-4.43655453180794},
    @{row=14; C=1464750; D=1294333; E=1391537; F=-4.99832735961768;  AAPC=-2.61146229525667; ICInf=-2.72616277238419; ICSup=-2.49662656902386},
    @{row=15; C=159070;  D=194605;  E=223558;  F=40.5406424844408;   AAPC=18.3373689865152;  ICInf=17.9614298700046;  ICSup=18.7145062083681},
    @{row=16; C=1305680; D=1099728; E=1167979; F=-10.5463053734453;  AAPC=-5.61942825253697; ICInf=-5.7393018067471;  ICSup=-5.49940225229924}
)

foreach ($r in $newData) {
    $row = $r.row
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 8).Value = $r.AAPC
    $ws.Cells.Item($row, 9).Value = $r.ICInf
    $ws.Cells.Item($row, 10).Value = $r.ICSup
}

# Drop the now-unused empty "delta2" column G; this shifts H->G, I->H, J->I,
# turning the old AAPC/IC_Inf/IC_Sup (H/I/J) layout into the new G/H/I one.
$ws.Columns.Item(7).Delete()

$ws.Range("M16").Select()
